$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "2025/2026, 2022/2023"
$ws.Range("G22").Value = "2025/2026, 2024/2025"
$ws.Range("G23").Value = "2025/2026, 2022/2023, 2023/2024"
$ws.Range("G24").Value = "2025/2026, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G27").Value = "2025/2026, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G28").Value = "2025/2026, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G31").Value = "2025/2026, 2022/2023"
$ws.Range("G50").Value = "2025/2026, 2024/2025"
$ws.Range("G51").Value = "2025/2026, 2022/2023, 2023/2024"
$ws.Range("G52").Value = "2025/2026, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G55").Value = "2025/2026, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G56").Value = "2025/2026, neveen.nashaat@med.asu.edu.eg"
